$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1089.3846
$ws.Range("I39").Value = 1647.8572
$ws.Range("J39").Value = 437.83334
$ws.Range("K39").Value = 4943.571599999999
$ws.Range("L39").Value = 1313.50002
$ws.Range("M39").Value = -4647.571599999999
$ws.Range("N39").Value = -1905.50002

$ws.Range("H100").Value = 12823272
$ws.Range("I100").Value = 33335414
$ws.Range("J100").Value = 3182.5
$ws.Range("K100").Value = 33335414
$ws.Range("L100").Value = 3182.5
$ws.Range("M100").Value = -33334873
$ws.Range("N100").Value = -4264.5

$ws.Range("H118").Value = 508.33334
$ws.Range("I118").Value = 450
$ws.Range("J118").Value = 800
$ws.Range("K118").Value = 1350
$ws.Range("L118").Value = 2400
$ws.Range("M118").Value = 307
$ws.Range("N118").Value = -5714

$ws.Range("H132").Value = 30170
$ws.Range("I132").Value = 34545.418
$ws.Range("K132").Value = 103636.254
$ws.Range("M132").Value = -101106.254

$ws.Range("H133").Value = 45507.617
$ws.Range("J133").Value = 45507.617
$ws.Range("L133").Value = 45507.617
$ws.Range("N133").Value = -55627.617

$ws.Range("H137").Value = 24391572
$ws.Range("I137").Value = 37038056
$ws.Range("J137").Value = 1927.3572
$ws.Range("K137").Value = 111114168
$ws.Range("L137").Value = 5782.071599999999
$ws.Range("M137").Value = -111111618
$ws.Range("N137").Value = -10882.0716

$ws.Range("H138").Value = 7502729.5
$ws.Range("I138").Value = 2157747.8
$ws.Range("J138").Value = 10872392
$ws.Range("K138").Value = 6473243.399999999
$ws.Range("L138").Value = 32617176
$ws.Range("M138").Value = -6468103.399999999
$ws.Range("N138").Value = -32627456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11078.5
$ws.Range("I28").Value = 2867.75
$ws.Range("J28").Value = 27500
$ws.Range("K28").Value = 2867.75
$ws.Range("L28").Value = 27500
$ws.Range("M28").Value = -2675.75
$ws.Range("N28").Value = -27884

$ws.Range("H32").Value = 25607.879
$ws.Range("I32").Value = 5260.6445
$ws.Range("J32").Value = 96040.62
$ws.Range("K32").Value = 5260.6445
$ws.Range("L32").Value = 96040.62
$ws.Range("M32").Value = -4973.6445
$ws.Range("N32").Value = -96614.62

$ws.Range("H61").Value = 2498.182
$ws.Range("I61").Value = 1771.7059
$ws.Range("K61").Value = 1771.7059
$ws.Range("M61").Value = -1559.7059

$ws.Range("H88").Value = 2007
$ws.Range("J88").Value = 2007
$ws.Range("L88").Value = 2007
$ws.Range("N88").Value = -2819

$ws.Range("H91").Value = 2007
$ws.Range("J91").Value = 2007
$ws.Range("L91").Value = 2007
$ws.Range("N91").Value = -4815

$ws.Range("H97").Value = 25750
$ws.Range("I97").Value = 100000
$ws.Range("K97").Value = 100000
$ws.Range("M97").Value = -99504

$ws.Range("H99").Value = 11078.5
$ws.Range("I99").Value = 2867.75
$ws.Range("J99").Value = 27500
$ws.Range("K99").Value = 2867.75
$ws.Range("L99").Value = 27500
$ws.Range("M99").Value = 127.25
$ws.Range("N99").Value = -33490

$ws.Range("H136").Value = 2498.182
$ws.Range("I136").Value = 1771.7059
$ws.Range("K136").Value = 5315.1177
$ws.Range("M136").Value = -2765.1177

$ws.Range("H139").Value = 46710
$ws.Range("J139").Value = 46710
$ws.Range("L139").Value = 46710
$ws.Range("N139").Value = -56990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 38755.668
$ws.Range("I26").Value = 26992.334
$ws.Range("J26").Value = 50519
$ws.Range("K26").Value = 26992.334
$ws.Range("L26").Value = 50519
$ws.Range("M26").Value = -26700.334
$ws.Range("N26").Value = -51103

$ws.Range("H59").Value = 43500
$ws.Range("J59").Value = 43500
$ws.Range("L59").Value = 43500
$ws.Range("N59").Value = -45194

$ws.Range("H134").Value = 3302.6365
$ws.Range("I134").Value = 2328.5217
$ws.Range("J134").Value = 5543.1
$ws.Range("K134").Value = 6985.5651
$ws.Range("L134").Value = 16629.3
$ws.Range("M134").Value = -4450.5651
$ws.Range("N134").Value = -21699.3

$ws.Range("H139").Value = 91260
$ws.Range("J139").Value = 91260
$ws.Range("L139").Value = 91260
$ws.Range("N139").Value = -101540

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 34921.5
$ws.Range("J88").Value = 34921.5
$ws.Range("L88").Value = 34921.5
$ws.Range("N88").Value = -35733.5

$ws.Range("H91").Value = 34921.5
$ws.Range("J91").Value = 34921.5
$ws.Range("L91").Value = 34921.5
$ws.Range("N91").Value = -37729.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5953773
$ws.Range("J131").Value = 7409050.5
$ws.Range("L131").Value = 22227151.5
$ws.Range("N131").Value = -22237231.5

$ws.Range("H136").Value = 2753.4443
$ws.Range("I136").Value = 1954
$ws.Range("J136").Value = 2853.375
$ws.Range("K136").Value = 5862
$ws.Range("L136").Value = 8560.125
$ws.Range("M136").Value = -762
$ws.Range("N136").Value = -18760.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1701.4286
$ws.Range("I97").Value = 1702.5
$ws.Range("J97").Value = 1700
$ws.Range("K97").Value = 1702.5
$ws.Range("L97").Value = 1700
$ws.Range("M97").Value = -1206.5
$ws.Range("N97").Value = -2692

$ws.Range("H102").Value = 2957.6428
$ws.Range("I102").Value = 2611.9443
$ws.Range("J102").Value = 3579.9
$ws.Range("K102").Value = 2611.9443
$ws.Range("L102").Value = 3579.9
$ws.Range("M102").Value = -989.9443000000001
$ws.Range("N102").Value = -6823.9

$ws.Range("H137").Value = 50226.668
$ws.Range("J137").Value = 50226.668
$ws.Range("L137").Value = 50226.668
$ws.Range("N137").Value = -60426.668

$ws.Range("H138").Value = 77999.664
$ws.Range("J138").Value = 77999.664
$ws.Range("L138").Value = 77999.664
$ws.Range("N138").Value = -88279.664

$ws.Range("H139").Value = 44996.5
$ws.Range("J139").Value = 44996.5
$ws.Range("L139").Value = 44996.5
$ws.Range("N139").Value = -55276.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2748.742
$ws.Range("I100").Value = 1956.3334
$ws.Range("J100").Value = 3072.9092
$ws.Range("K100").Value = 1956.3334
$ws.Range("L100").Value = 3072.9092
$ws.Range("M100").Value = -1415.3334
$ws.Range("N100").Value = -4154.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 900
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -2482
